$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("csv_1")

# Fill in values 1-9 across columns B through J for rows 2, 3, 4
$values = @(1,2,3,4,5,6,7,8,9)
for ($row = 2; $row -le 4; $row++) {
    for ($col = 2; $col -le 10; $col++) {
        $ws1.Cells.Item($row, $col).Value = $values[$col - 2]
    }
}

# Update the active selection on sheet1 to G4
$ws1.Range("G4").Select()

$wb.Save()
